$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Componentes Conexos" column with new percolation-by-nodes values
$ws.Range("B2").Value = 9830
$ws.Range("B3").Value = 9268
$ws.Range("B4").Value = 8331
$ws.Range("B5").Value = 7149
$ws.Range("B6").Value = 5730
$ws.Range("B7").Value = 4225
$ws.Range("B8").Value = 3052
$ws.Range("B9").Value = 1971
$ws.Range("B10").Value = 962
$ws.Range("B11").Value = 1

# Remove the now-obsolete extra blocks of rows (rows 12 through 33)
$ws.Range("A12:B33").EntireRow.Delete() | Out-Null
